$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Premier League "rodada 7" standings table (rows 2-21, columns A-I)
# A=clube, B=pontos, C=partidas jogadas, D=vitorias, E=empates, F=derrotas,
# G=total de gols, H=total de gols sofridos, I=saldo de gols
$data = @(
    @{r=2;  A="Liverpool";          B=18; C=7; D=6; E=0; F=1; G=13; H=2;  I=11}
    @{r=3;  A="Arsenal";            B=17; C=7; D=5; E=2; F=0; G=15; H=6;  I=9}
    @{r=4;  A="City";               B=17; C=7; D=5; E=2; F=0; G=17; H=8;  I=9}
    @{r=5;  A="Chelsea";            B=14; C=7; D=4; E=2; F=1; G=16; H=8;  I=8}
    @{r=6;  A="Aston Villa";        B=14; C=7; D=4; E=2; F=1; G=12; H=9;  I=3}
    @{r=7;  A="Brighton";           B=12; C=7; D=3; E=3; F=1; G=13; H=10; I=3}
    @{r=8;  A="Newcastle";          B=12; C=7; D=3; E=3; F=1; G=8;  H=7;  I=1}
    @{r=9;  A="Fulham";             B=11; C=7; D=3; E=2; F=2; G=10; H=8;  I=2}
    @{r=10; A="Tottenham";          B=10; C=7; D=3; E=1; F=3; G=14; H=8;  I=6}
    @{r=11; A="Brentford";          B=10; C=7; D=3; E=1; F=3; G=13; H=13; I=0}
    @{r=12; A="Nottingham Forest";  B=10; C=7; D=2; E=4; F=1; G=7;  H=6;  I=1}
    @{r=13; A="West Ham";           B=8;  C=7; D=2; E=2; F=3; G=10; H=11; I=-1}
    @{r=14; A="Bournemouth";        B=8;  C=7; D=2; E=2; F=3; G=8;  H=10; I=-2}
    @{r=15; A="Manchester United";  B=8;  C=7; D=2; E=2; F=3; G=5;  H=8;  I=-3}
    @{r=16; A="Leicester City";     B=6;  C=7; D=1; E=3; F=3; G=9;  H=12; I=-3}
    @{r=17; A="Everton";            B=5;  C=7; D=1; E=2; F=4; G=7;  H=15; I=-8}
    @{r=18; A="Ipswich Town";       B=4;  C=7; D=0; E=4; F=3; G=6;  H=14; I=-8}
    @{r=19; A="Crystal Palace";     B=3;  C=7; D=0; E=3; F=4; G=5;  H=10; I=-5}
    @{r=20; A="Southampton";        B=1;  C=7; D=0; E=1; F=6; G=4;  H=15; I=-11}
    @{r=21; A="Wolves";             B=1;  C=7; D=0; E=1; F=6; G=9;  H=21; I=-12}
)

foreach ($row in $data) {
    $r = $row.r
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
}
